$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 11.98327633333333
$ws.Range("H2").Value = 35.949829
$ws.Range("I2").Value = 0.03345300399843466
$ws.Range("J2").Value = 0.03345300399843466
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 30.99161333333333
$ws.Range("N2").Value = 92.97484
$ws.Range("O2").Value = 0.3599121977633812
$ws.Range("P2").Value = 0.3599121977633811
$ws.Range("Q2").Value = 371.3810665891511
$ws.Range("R2").Value = 3342.42959930236
$ws.Range("S2").Value = 0.0120401441908638
$ws.Range("T2").Value = 0.01204014419086379

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 11.98327633333333
$ws.Range("H3").Value = 35.949829
$ws.Range("I3").Value = 0.03345300399843466
$ws.Range("J3").Value = 0.03345300399843466
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 29.913269
$ws.Range("N3").Value = 89.739807
$ws.Range("O3").Value = 0.3473891556493311
$ws.Range("P3").Value = 0.3473891556493311
$ws.Range("Q3").Value = 358.4589684603337
$ws.Range("R3").Value = 3226.130716143003
$ws.Range("S3").Value = 0.01162121081294991
$ws.Range("T3").Value = 0.01162121081294991

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 11.98327633333333
$ws.Range("H4").Value = 35.949829
$ws.Range("I4").Value = 0.03345300399843466
$ws.Range("J4").Value = 0.03345300399843466
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 19.150218
$ws.Range("N4").Value = 57.450654
$ws.Range("O4").Value = 0.2223955550134164
$ws.Range("P4").Value = 0.2223955550134163
$ws.Range("Q4").Value = 229.482354137574
$ws.Range("R4").Value = 2065.341187238166
$ws.Range("S4").Value = 0.007439799391097913
$ws.Range("T4").Value = 0.007439799391097911

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 11.98327633333333
$ws.Range("H5").Value = 35.949829
$ws.Range("I5").Value = 0.03345300399843466
$ws.Range("J5").Value = 0.03345300399843466
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 6.053716000000001
$ws.Range("N5").Value = 18.161148
$ws.Range("O5").Value = 0.07030309157387134
$ws.Range("P5").Value = 0.07030309157387132
$ws.Range("Q5").Value = 72.54335167152135
$ws.Range("R5").Value = 652.890165043692
$ws.Range("S5").Value = 0.002351849603523036
$ws.Range("T5").Value = 0.002351849603523035

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 332.3726806666667
$ws.Range("H6").Value = 997.1180420000001
$ws.Range("I6").Value = 0.9278651602470024
$ws.Range("J6").Value = 0.9278651602470025
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 30.99161333333333
$ws.Range("N6").Value = 92.97484
$ws.Range("O6").Value = 0.3599121977633812
$ws.Range("P6").Value = 0.3599121977633811
$ws.Range("Q6").Value = 10300.76560178481
$ws.Range("R6").Value = 92706.89041606328
$ws.Range("S6").Value = 0.3339499890525706
$ws.Range("T6").Value = 0.3339499890525705

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 332.3726806666667
$ws.Range("H7").Value = 997.1180420000001
$ws.Range("I7").Value = 0.9278651602470024
$ws.Range("J7").Value = 0.9278651602470025
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 29.913269
$ws.Range("N7").Value = 89.739807
$ws.Range("O7").Value = 0.3473891556493311
$ws.Range("P7").Value = 0.3473891556493311
$ws.Range("Q7").Value = 9942.3534050331
$ws.Range("R7").Value = 89481.18064529789
$ws.Range("S7").Value = 0.3223302945746375
$ws.Range("T7").Value = 0.3223302945746375

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 332.3726806666667
$ws.Range("H8").Value = 997.1180420000001
$ws.Range("I8").Value = 0.9278651602470024
$ws.Range("J8").Value = 0.9278651602470025
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 19.150218
$ws.Range("N8").Value = 57.450654
$ws.Range("O8").Value = 0.2223955550134164
$ws.Range("P8").Value = 0.2223955550134163
$ws.Range("Q8").Value = 6365.009292011052
$ws.Range("R8").Value = 57285.08362809947
$ws.Range("S8").Value = 0.2063530872907446
$ws.Range("T8").Value = 0.2063530872907446

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 332.3726806666667
$ws.Range("H9").Value = 997.1180420000001
$ws.Range("I9").Value = 0.9278651602470024
$ws.Range("J9").Value = 0.9278651602470025
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 6.053716000000001
$ws.Range("N9").Value = 18.161148
$ws.Range("O9").Value = 0.07030309157387134
$ws.Range("P9").Value = 0.07030309157387132
$ws.Range("Q9").Value = 2012.089814914691
$ws.Range("R9").Value = 18108.80833423222
$ws.Range("S9").Value = 0.06523178932904981
$ws.Range("T9").Value = 0.0652317893290498

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.08615933333333332
$ws.Range("H10").Value = 0.258478
$ws.Range("I10").Value = 0.0002405259164795302
$ws.Range("J10").Value = 0.0002405259164795302
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 30.99161333333333
$ws.Range("N10").Value = 92.97484
$ws.Range("O10").Value = 0.3599121977633812
$ws.Range("P10").Value = 0.3599121977633811
$ws.Range("Q10").Value = 2.670216743724444
$ws.Range("R10").Value = 24.03195069352
$ws.Range("S10").Value = 0.00008656821121919917
$ws.Range("T10").Value = 0.00008656821121919916

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.08615933333333332
$ws.Range("H11").Value = 0.258478
$ws.Range("I11").Value = 0.0002405259164795302
$ws.Range("J11").Value = 0.0002405259164795302
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 29.913269
$ws.Range("N11").Value = 89.739807
$ws.Range("O11").Value = 0.3473891556493311
$ws.Range("P11").Value = 0.3473891556493311
$ws.Range("Q11").Value = 2.577307314860666
$ws.Range("R11").Value = 23.195765833746
$ws.Range("S11").Value = 0.00008355609503760552
$ws.Range("T11").Value = 0.00008355609503760552

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.08615933333333332
$ws.Range("H12").Value = 0.258478
$ws.Range("I12").Value = 0.0002405259164795302
$ws.Range("J12").Value = 0.0002405259164795302
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 19.150218
$ws.Range("N12").Value = 57.450654
$ws.Range("O12").Value = 0.2223955550134164
$ws.Range("P12").Value = 0.2223955550134163
$ws.Range("Q12").Value = 1.649970016068
$ws.Range("R12").Value = 14.849730144612
$ws.Range("S12").Value = 0.00005349189469057575
$ws.Range("T12").Value = 0.00005349189469057574

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.08615933333333332
$ws.Range("H13").Value = 0.258478
$ws.Range("I13").Value = 0.0002405259164795302
$ws.Range("J13").Value = 0.0002405259164795302
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 6.053716000000001
$ws.Range("N13").Value = 18.161148
$ws.Range("O13").Value = 0.07030309157387134
$ws.Range("P13").Value = 0.07030309157387132
$ws.Range("Q13").Value = 0.5215841347493333
$ws.Range("R13").Value = 4.694257212744
$ws.Range("S13").Value = 0.00001690971553214974
$ws.Range("T13").Value = 0.00001690971553214974

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 13.77014866666667
$ws.Range("H14").Value = 41.310446
$ws.Range("I14").Value = 0.03844130983808348
$ws.Range("J14").Value = 0.03844130983808348
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 30.99161333333333
$ws.Range("N14").Value = 92.97484
$ws.Range("O14").Value = 0.3599121977633812
$ws.Range("P14").Value = 0.3599121977633811
$ws.Range("Q14").Value = 426.7591230198489
$ws.Range("R14").Value = 3840.83210717864
$ws.Range("S14").Value = 0.01383549630872771
$ws.Range("T14").Value = 0.01383549630872771

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 13.77014866666667
$ws.Range("H15").Value = 41.310446
$ws.Range("I15").Value = 0.03844130983808348
$ws.Range("J15").Value = 0.03844130983808348
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 29.913269
$ws.Range("N15").Value = 89.739807
$ws.Range("O15").Value = 0.3473891556493311
$ws.Range("P15").Value = 0.3473891556493311
$ws.Range("Q15").Value = 411.9101612359913
$ws.Range("R15").Value = 3707.191451123922
$ws.Range("S15").Value = 0.01335409416670615
$ws.Range("T15").Value = 0.01335409416670615

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 13.77014866666667
$ws.Range("H16").Value = 41.310446
$ws.Range("I16").Value = 0.03844130983808348
$ws.Range("J16").Value = 0.03844130983808348
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 19.150218
$ws.Range("N16").Value = 57.450654
$ws.Range("O16").Value = 0.2223955550134164
$ws.Range("P16").Value = 0.2223955550134163
$ws.Range("Q16").Value = 263.7013488590759
$ws.Range("R16").Value = 2373.312139731684
$ws.Range("S16").Value = 0.00854917643688328
$ws.Range("T16").Value = 0.008549176436883278

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 13.77014866666667
$ws.Range("H17").Value = 41.310446
$ws.Range("I17").Value = 0.03844130983808348
$ws.Range("J17").Value = 0.03844130983808348
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 6.053716000000001
$ws.Range("N17").Value = 18.161148
$ws.Range("O17").Value = 0.07030309157387134
$ws.Range("P17").Value = 0.07030309157387132
$ws.Range("Q17").Value = 83.36056930577867
$ws.Range("R17").Value = 750.245123752008
$ws.Range("S17").Value = 0.002702542925766344
$ws.Range("T17").Value = 0.002702542925766344
